$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_coa")

$data = @(
    ,@(2, "Keperluan Komputer", "Keperluan Komputer", "COA untuk Keperluan Komputer")
    ,@(3, "M. HW", "Maintenance Hardware", "COA untuk Maintenance Hardware")
    ,@(4, "Test COA", "Hype COA", "Def COA")
    ,@(5, "Coa 1", "Coa 3", "Coa 2")
    ,@(6, "Coa 2", "Coa 2", "Coa 2")
    ,@(7, "Coa 3", "Coa 3", "Coa 3")
    ,@(8, "Coa 4", "Coa 4", "Coa 4")
    ,@(9, "Coa 5", "Coa 5", "Coa 5")
    ,@(10, "Gedung", "Gedung Hyperion", "COA untuk Gedung")
    ,@(11, "ASET ROA TANAH", "ROA", "COA untuk Roa Tanah")
    ,@(12, "Pemeliharaan Gedung & Perabotan", "M.Gedung&Perabotan", "COA untuk Maintenance Gedung & Perabotan")
    ,@(13, "SW", "Software", "COA untuk Software")
    ,@(14, "HW", "Hardware", "COA untuk Hardware")
    ,@(15, "M. SW", "Maintenance Software", "COA untuk Maintenance Software")
    ,@(16, "Consultant", "Consultant", "COA untuk Consultant")
    ,@(17, "Pemanfaatan IT", "Pemanfaatan IT", "COA untuk pemanfaatan IT")
    ,@(18, "Sewa Gedung", "Sewa Gedung", "COA untuk Sewa Gedung")
    ,@(19, "Keperluan Kantor Lainnya", "Miscellaneous", "COA untuk keperluan kantor lainnya")
    ,@(20, "Biaya Pengelola Pendukung Operasional", "Biaya Pengelola Pendukung Operator", "COA Pendukung Operator")
    ,@(21, "Promosi", "Promosi", "COA untuk Promosi")
    ,@(22, "Biaya Pengelola Pendukung Proses", "Pengelola Pendukung Proses", "COA untuk Biaya Pengelola Pendukung Proses")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
